# Generate Report for Handoff
#
# The "48bc4e3f-4ec2-44de-a3ab-d3ced75e777a" file has dropped out of the
# report (its row is removed from every sheet), and the
# "01a7335d-4021-43b2-ba3d-305a29b44724" file moved from "handed back" to
# "ready for handoff" with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

function Remove-RowHyperlinks($ws, $rowNum) {
    # Deleting items from a live COM collection while enumerating it is
    # unreliable here, so repeatedly rescan-and-delete-one until none of
    # the target row's hyperlinks remain.
    $found = $true
    while ($found) {
        $found = $false
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Row -eq $rowNum) {
                $hl.Delete()
                $found = $true
                break
            }
        }
    }
}

# --- Sheet "Overview" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-03-20 04:47:27"

Remove-RowHyperlinks $ws1 3
$ws1.Rows.Item(3).Delete()

# --- Sheet "zh-cn" ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-20 04:47:19"

Remove-RowHyperlinks $ws2 3
$ws2.Rows.Item(3).Delete()

# --- Sheet "de-de" ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-20 04:47:27"

Remove-RowHyperlinks $ws3 3
$ws3.Rows.Item(3).Delete()
